$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 (NB) with refreshed measurement values
$ws.Range("B5").Value = 2511.0
$ws.Range("C5").Value = 3594.0
$ws.Range("D5").Value = 32954.0
$ws.Range("E5").Value = 2129.0
$ws.Range("F5").Value = 0.4113022113022113
$ws.Range("G5").Value = 0.5411637931034483
$ws.Range("H5").Value = 0.9016635657217905
$ws.Range("I5").Value = 0.46738017682643085
$ws.Range("J5").Value = 0.8610517626493154
$ws.Range("K5").Value = 0.3891878394859579
$ws.Range("L5").Value = 0.6985325155847683
$ws.Range("M5").Value = 0.8335700182187566
